$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 295-296, pushing existing data (old rows 295-409)
# down to become rows 297-411.
$ws.Rows("295:296").Insert()

# Populate the newly inserted row 295 with new weekly data.
$ws.Cells.Item(295, 1).Value = 7
$ws.Cells.Item(295, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(295, 3).Value = "Ñuble"
$ws.Cells.Item(295, 4).Value = 45027
$ws.Cells.Item(295, 5).Value = 16
$ws.Cells.Item(295, 6).Value = 100114013
$ws.Cells.Item(295, 7).Value = "Zanahoria"
$ws.Cells.Item(295, 8).Value = "Sin especificar"
$ws.Cells.Item(295, 9).Value = "Primera"
$ws.Cells.Item(295, 10).Value = 150
$ws.Cells.Item(295, 11).Value = 7000
$ws.Cells.Item(295, 12).Value = 7000
$ws.Cells.Item(295, 13).Value = 7000
$ws.Cells.Item(295, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(295, 15).Value = "Región de Ñuble"
$ws.Cells.Item(295, 16).Value = 350
$ws.Cells.Item(295, 17).Value = 20
$ws.Cells.Item(295, 18).Value = "Hortaliza"

# Populate the newly inserted row 296 with new weekly data.
$ws.Cells.Item(296, 1).Value = 7
$ws.Cells.Item(296, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(296, 3).Value = "Ñuble"
$ws.Cells.Item(296, 4).Value = 45027
$ws.Cells.Item(296, 5).Value = 16
$ws.Cells.Item(296, 6).Value = 100114013
$ws.Cells.Item(296, 7).Value = "Zanahoria"
$ws.Cells.Item(296, 8).Value = "Sin especificar"
$ws.Cells.Item(296, 9).Value = "Segunda"
$ws.Cells.Item(296, 10).Value = 200
$ws.Cells.Item(296, 11).Value = 6000
$ws.Cells.Item(296, 12).Value = 6000
$ws.Cells.Item(296, 13).Value = 6000
$ws.Cells.Item(296, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(296, 15).Value = "Región de Ñuble"
$ws.Cells.Item(296, 16).Value = 300
$ws.Cells.Item(296, 17).Value = 20
$ws.Cells.Item(296, 18).Value = "Hortaliza"
